$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "IDN"
$ws.Range("B8").Value = "räksmörgås.josefsson.org"
$ws.Range("C8").Value = "räksmörgås.josefsson.org"
$ws.Range("D8").Value = "Test"

$ws.Range("D8").Select()
